$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4675
$ws.Range("J10").Value = 4675
$ws.Range("L10").Value = 4675
$ws.Range("N10").Value = -5261
$ws.Range("H11").Value = 453.14285
$ws.Range("I11").Value = 453.14285
$ws.Range("K11").Value = 453.14285
$ws.Range("M11").Value = -313.14285
$ws.Range("H49").Value = 1225
$ws.Range("I49").Value = 600
$ws.Range("J49").Value = 1850
$ws.Range("K49").Value = 1800
$ws.Range("L49").Value = 5550
$ws.Range("M49").Value = -1664
$ws.Range("N49").Value = -5822
$ws.Range("H96").Value = 588970.4399999999
$ws.Range("I96").Value = 769646.9
$ws.Range("J96").Value = 1772
$ws.Range("K96").Value = 2308940.7
$ws.Range("L96").Value = 5316
$ws.Range("M96").Value = -2307567.7
$ws.Range("N96").Value = -8062
$ws.Range("H98").Value = 7655.304
$ws.Range("I98").Value = 10021.883
$ws.Range("K98").Value = 10021.883
$ws.Range("M98").Value = -8523.883
$ws.Range("H122").Value = 7655.304
$ws.Range("I122").Value = 10021.883
$ws.Range("K122").Value = 30065.649
$ws.Range("M122").Value = -27615.649
$ws.Range("H132").Value = 2143.1785
$ws.Range("I132").Value = 1713.5
$ws.Range("J132").Value = 3217.375
$ws.Range("K132").Value = 5140.5
$ws.Range("L132").Value = 9652.125
$ws.Range("M132").Value = -2610.5
$ws.Range("N132").Value = -14712.125
$ws.Range("H138").Value = 1101.5588
$ws.Range("I138").Value = 1032.6666
$ws.Range("J138").Value = 3375
$ws.Range("K138").Value = 3097.9998
$ws.Range("L138").Value = 10125
$ws.Range("M138").Value = 2042.0002
$ws.Range("N138").Value = -20405

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 6126
$ws.Range("I22").Value = 5534.4
$ws.Range("J22").Value = 15000
$ws.Range("K22").Value = 5534.4
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = -5235.4
$ws.Range("N22").Value = -15598
$ws.Range("H88").Value = 1477.6923
$ws.Range("I88").Value = 1297
$ws.Range("J88").Value = 1558
$ws.Range("K88").Value = 1297
$ws.Range("L88").Value = 1558
$ws.Range("M88").Value = -891
$ws.Range("N88").Value = -2370
$ws.Range("H91").Value = 1477.6923
$ws.Range("I91").Value = 1297
$ws.Range("J91").Value = 1558
$ws.Range("K91").Value = 1297
$ws.Range("L91").Value = 1558
$ws.Range("M91").Value = 107
$ws.Range("N91").Value = -4366
$ws.Range("H122").Value = 18520408
$ws.Range("I122").Value = 18520408
$ws.Range("K122").Value = 55561224
$ws.Range("M122").Value = -55558774
$ws.Range("H132").Value = 1622.3934
$ws.Range("I132").Value = 1651.0186
$ws.Range("K132").Value = 4953.0558
$ws.Range("M132").Value = -2423.0558

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 40.4
$ws.Range("I12").Value = 38
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 38
$ws.Range("L12").Value = 50
$ws.Range("M12").Value = 130
$ws.Range("N12").Value = -386
$ws.Range("H134").Value = 2227.7778
$ws.Range("I134").Value = 2102.1738
$ws.Range("K134").Value = 6306.5214
$ws.Range("M134").Value = -3771.5214

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2344.2222
$ws.Range("I31").Value = 1119.8
$ws.Range("K31").Value = 1119.8
$ws.Range("M31").Value = -824.8
$ws.Range("H34").Value = 2344.2222
$ws.Range("I34").Value = 1119.8
$ws.Range("K34").Value = 1119.8
$ws.Range("M34").Value = -917.8
$ws.Range("H58").Value = 3407.6428
$ws.Range("I58").Value = 2839.3
$ws.Range("K58").Value = 2839.3
$ws.Range("M58").Value = -2636.3
$ws.Range("H136").Value = 3407.6428
$ws.Range("I136").Value = 2839.3
$ws.Range("K136").Value = 8517.900000000001
$ws.Range("M136").Value = -5967.900000000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 120.333336
$ws.Range("J2").Value = 92
$ws.Range("L2").Value = 552
$ws.Range("N2").Value = -778
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H116").Value = 800
$ws.Range("I116").Value = 800
$ws.Range("K116").Value = 2400
$ws.Range("M116").Value = 1042
$ws.Range("H138").Value = 5255.6665
$ws.Range("I138").Value = 3435
$ws.Range("K138").Value = 10305
$ws.Range("M138").Value = -5165

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H70").Value = 111116310
$ws.Range("J70").Value = 7600
$ws.Range("L70").Value = 7600
$ws.Range("N70").Value = -8140
$ws.Range("H73").Value = 111116310
$ws.Range("J73").Value = 7600
$ws.Range("L73").Value = 7600
$ws.Range("N73").Value = -9472
$ws.Range("H97").Value = 14351.182
$ws.Range("J97").Value = 30596.8
$ws.Range("L97").Value = 30596.8
$ws.Range("N97").Value = -31588.8
$ws.Range("H113").Value = 7131.96
$ws.Range("I113").Value = 3733.25
$ws.Range("J113").Value = 10269.23
$ws.Range("K113").Value = 3733.25
$ws.Range("L113").Value = 10269.23
$ws.Range("M113").Value = -1563.25
$ws.Range("N113").Value = -14609.23
$ws.Range("H132").Value = 2845.56
$ws.Range("I132").Value = 2732.6
$ws.Range("J132").Value = 3297.4
$ws.Range("K132").Value = 8197.799999999999
$ws.Range("L132").Value = 9892.200000000001
$ws.Range("M132").Value = -5667.799999999999
$ws.Range("N132").Value = -14952.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 820.625
$ws.Range("I16").Value = 650.75
$ws.Range("J16").Value = 990.5
$ws.Range("K16").Value = 650.75
$ws.Range("L16").Value = 990.5
$ws.Range("M16").Value = -480.75
$ws.Range("N16").Value = -1330.5
$ws.Range("H93").Value = 3938.577
$ws.Range("I93").Value = 1413
$ws.Range("K93").Value = 1413
$ws.Range("M93").Value = -165

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2578.5386
$ws.Range("I96").Value = 1315.125
$ws.Range("K96").Value = 1315.125
$ws.Range("M96").Value = 57.875
